# Full dataset & dara corr.
#
# The data that used to live in columns B:E of rows 43-69 needs to move down
# by exactly one row (into rows 44-70); row 43's B:E become blank again
# (that statistic got dropped from the top of the table). Column A (the
# running index numbers 41, 42, 43, ...) is untouched throughout.
#
# We shift row-by-row, bottom to top, so a row is always fully read (via
# Copy) before it is overwritten. Using Copy/PasteSpecial(values) instead of
# a plain .Value assignment keeps numeric-looking text (e.g. "12581") stored
# as text/shared-string, matching the source data instead of being
# auto-coerced into a real number by the usual COM .Value setter. Each
# destination range is cleared first so that rows whose source is partially
# blank (only column B populated) don't leave stale values behind in C:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 70; $r -ge 44; $r--) {
    $src = $r - 1
    $destRange = $ws.Range("B$r" + ":E$r")
    $destRange.ClearContents()
    $ws.Range("B$src" + ":E$src").Copy()
    $destRange.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false

# Row 43's B:E is now empty (its old content moved down into row 44).
$ws.Range("B43:E43").ClearContents()

# Match the author's resulting view state: scrolled so row 39 is at the top,
# with G44 as the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G44").Select()
